$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 10685.714
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2700
$ws.Range("I63").Value = 3200
$ws.Range("J63").Value = 1450
$ws.Range("K63").Value = 3200
$ws.Range("L63").Value = 1450
$ws.Range("M63").Value = -2514
$ws.Range("N63").Value = -2822
$ws.Range("H66").Value = 2700
$ws.Range("I66").Value = 3200
$ws.Range("J66").Value = 1450
$ws.Range("K66").Value = 16000
$ws.Range("L66").Value = 7250
$ws.Range("M66").Value = -12568
$ws.Range("N66").Value = -14114
$ws.Range("H74").Value = 307817.38
$ws.Range("I74").Value = 436337.7
$ws.Range("J74").Value = 80435.234
$ws.Range("K74").Value = 436337.7
$ws.Range("L74").Value = 80435.234
$ws.Range("M74").Value = -435463.7
$ws.Range("N74").Value = -82183.234
$ws.Range("H77").Value = 307817.38
$ws.Range("I77").Value = 436337.7
$ws.Range("J77").Value = 80435.234
$ws.Range("K77").Value = 2181688.5
$ws.Range("L77").Value = 402176.17
$ws.Range("M77").Value = -2177320.5
$ws.Range("N77").Value = -410912.17
$ws.Range("H132").Value = 27857.61
$ws.Range("I132").Value = 40156.297
$ws.Range("J132").Value = 4138.7144
$ws.Range("K132").Value = 120468.891
$ws.Range("L132").Value = 12416.1432
$ws.Range("M132").Value = -117938.891
$ws.Range("N132").Value = -17476.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 889.825
$ws.Range("I94").Value = 653.2963
$ws.Range("J94").Value = 1381.0769
$ws.Range("K94").Value = 653.2963
$ws.Range("L94").Value = 1381.0769
$ws.Range("M94").Value = -202.2963
$ws.Range("N94").Value = -2283.0769
$ws.Range("H131").Value = 25726.666
$ws.Range("J131").Value = 25726.666
$ws.Range("L131").Value = 25726.666
$ws.Range("N131").Value = -35806.666
$ws.Range("H134").Value = 3880.9395
$ws.Range("I134").Value = 3720.4783
$ws.Range("J134").Value = 4250
$ws.Range("K134").Value = 11161.4349
$ws.Range("L134").Value = 12750
$ws.Range("M134").Value = -8626.4349
$ws.Range("N134").Value = -17820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 57.77778
$ws.Range("I7").Value = 37.5
$ws.Range("J7").Value = 98.333336
$ws.Range("K7").Value = 37.5
$ws.Range("L7").Value = 98.333336
$ws.Range("M7").Value = 75.5
$ws.Range("N7").Value = -324.333336
$ws.Range("H22").Value = 210.03448
$ws.Range("I22").Value = 105.5625
$ws.Range("K22").Value = 105.5625
$ws.Range("M22").Value = 244.4375
$ws.Range("H38").Value = 4600
$ws.Range("I38").Value = 2000
$ws.Range("J38").Value = 6333.3335
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 6333.3335
$ws.Range("M38").Value = -1623
$ws.Range("N38").Value = -7087.3335
$ws.Range("H46").Value = 4600
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 6333.3335
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 6333.3335
$ws.Range("M46").Value = -1789
$ws.Range("N46").Value = -6755.3335
$ws.Range("H58").Value = 3479.7778
$ws.Range("I58").Value = 4071.5806
$ws.Range("K58").Value = 4071.5806
$ws.Range("M58").Value = -3868.5806
$ws.Range("H107").Value = 208.92857
$ws.Range("I107").Value = 158.43478
$ws.Range("J107").Value = 441.2
$ws.Range("K107").Value = 158.43478
$ws.Range("L107").Value = 441.2
$ws.Range("M107").Value = 1761.56522
$ws.Range("N107").Value = -4281.2
$ws.Range("H122").Value = 1135.2354
$ws.Range("I122").Value = 1042.6364
$ws.Range("J122").Value = 1305
$ws.Range("K122").Value = 3127.9092
$ws.Range("L122").Value = 3915
$ws.Range("M122").Value = -677.9092000000001
$ws.Range("N122").Value = -8815
$ws.Range("H134").Value = 1329.6938
$ws.Range("I134").Value = 871.7568
$ws.Range("K134").Value = 2615.2704
$ws.Range("M134").Value = -80.27039999999988
$ws.Range("H136").Value = 3479.7778
$ws.Range("I136").Value = 4071.5806
$ws.Range("K136").Value = 12214.7418
$ws.Range("M136").Value = -9664.7418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 6731.933
$ws.Range("I2").Value = 14311.286
$ws.Range("K2").Value = 85867.716
$ws.Range("M2").Value = -85754.716
$ws.Range("H38").Value = 126.25
$ws.Range("I38").Value = 87.14286
$ws.Range("K38").Value = 261.42858
$ws.Range("M38").Value = 85.57141999999999
$ws.Range("H98").Value = 10398.5
$ws.Range("J98").Value = 12398.2
$ws.Range("L98").Value = 37194.60000000001
$ws.Range("N98").Value = -40190.60000000001
$ws.Range("H102").Value = 5650
$ws.Range("I102").Value = 3428.5715
$ws.Range("J102").Value = 8760
$ws.Range("K102").Value = 10285.7145
$ws.Range("L102").Value = 26280
$ws.Range("M102").Value = -7851.7145
$ws.Range("N102").Value = -31148
$ws.Range("H108").Value = 2661.9092
$ws.Range("I108").Value = 2586.7778
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 7760.3334
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = -4880.3334
$ws.Range("N108").Value = -14760
$ws.Range("H113").Value = 9091453
$ws.Range("I113").Value = 12500530
$ws.Range("J113").Value = 581.3333
$ws.Range("K113").Value = 37501590
$ws.Range("L113").Value = 1743.9999
$ws.Range("M113").Value = -37499420
$ws.Range("N113").Value = -6083.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3350.5972
$ws.Range("I80").Value = 3576.228
$ws.Range("J80").Value = 2493.2
$ws.Range("K80").Value = 3576.228
$ws.Range("L80").Value = 2493.2
$ws.Range("M80").Value = -2578.228
$ws.Range("N80").Value = -4489.2
$ws.Range("H83").Value = 3350.5972
$ws.Range("I83").Value = 3576.228
$ws.Range("J83").Value = 2493.2
$ws.Range("K83").Value = 17881.14
$ws.Range("L83").Value = 12466
$ws.Range("M83").Value = -12889.14
$ws.Range("N83").Value = -22450
$ws.Range("H113").Value = 1469.5
$ws.Range("I113").Value = 961
$ws.Range("J113").Value = 2268.5715
$ws.Range("K113").Value = 961
$ws.Range("L113").Value = 2268.5715
$ws.Range("M113").Value = 1209
$ws.Range("N113").Value = -6608.5715
$ws.Range("H132").Value = 3833
$ws.Range("I132").Value = 3637.077
$ws.Range("J132").Value = 4224.846
$ws.Range("K132").Value = 10911.231
$ws.Range("L132").Value = 12674.538
$ws.Range("M132").Value = -8381.231
$ws.Range("N132").Value = -17734.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 637.75
$ws.Range("I107").Value = 717
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 2151
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -231
$ws.Range("N107").Value = -5040
$ws.Range("H122").Value = 27028268
$ws.Range("I122").Value = 30304210
$ws.Range("J122").Value = 1740
$ws.Range("K122").Value = 90912630
$ws.Range("L122").Value = 5220
$ws.Range("M122").Value = -90910180
$ws.Range("N122").Value = -10120
$ws.Range("H136").Value = 21513086
$ws.Range("I136").Value = 31283198
$ws.Range("J136").Value = 670180.3
$ws.Range("K136").Value = 93849594
$ws.Range("L136").Value = 2010540.9
$ws.Range("M136").Value = -93847044
$ws.Range("N136").Value = -2015640.9
